$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C215:C252").Value = 7310
